$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -6
$ws.Range("F11").Value = 5
$ws.Range("F13").Value = -7
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 3
$ws.Range("F21").Value = -3
$ws.Range("F26").Value = -3
$ws.Range("F28").Value = 0
$ws.Range("F36").Value = -7
$ws.Range("F41").Value = -1
$ws.Range("F42").Value = 0
$ws.Range("F45").Value = 4
$ws.Range("F46").Value = -1
